$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 <- values previously in row 15
$ws.Range("A13").Value = 130864513
$ws.Range("B13").Value = 78255
$ws.Range("E13").Value = 228579
$ws.Range("F13").Value = "Liten svartspik"
$ws.Range("G13").Value = "Chaenothecopsis nana"
$ws.Range("H13").Value = "Tibell"
$ws.Range("Q13").Value = 445962
$ws.Range("R13").Value = 7031114
$ws.Range("S13").Value = 4
$ws.Range("Z13").Value = "10:51"
$ws.Range("AB13").Value = "10:51"

# Row 14 <- values previously in row 13
$ws.Range("A14").Value = 130864511
$ws.Range("B14").Value = 83223
$ws.Range("E14").Value = 6440
$ws.Range("F14").Value = "Vitgrynig nållav"
$ws.Range("G14").Value = "Chaenotheca subroscida"
$ws.Range("H14").Value = "(Eitner) Zahlbr."
$ws.Range("Q14").Value = 445985
$ws.Range("R14").Value = 7031157
$ws.Range("S14").Value = 3
$ws.Range("Z14").Value = "10:40"
$ws.Range("AB14").Value = "10:40"

# Row 15 <- values previously in row 14
$ws.Range("A15").Value = 130864514
$ws.Range("B15").Value = 83223
$ws.Range("E15").Value = 6440
$ws.Range("F15").Value = "Vitgrynig nållav"
$ws.Range("G15").Value = "Chaenotheca subroscida"
$ws.Range("H15").Value = "(Eitner) Zahlbr."
$ws.Range("Q15").Value = 445944
$ws.Range("R15").Value = 7031101
$ws.Range("S15").Value = 4
$ws.Range("Z15").Value = "11:03"
$ws.Range("AB15").Value = "11:03"
